$d = $word.ActiveDocument
$sel = $word.Selection

# Move to the very end of the document (past the final bookmarkEnd),
# so new paragraphs are appended after the "Students" heading.
$sel.EndKey(6)

# Phase 1: type all of the new paragraphs as plain text, remembering the
# start offset of each paragraph and how many leading characters (if any)
# need to end up bold. Bold formatting is applied afterwards (phase 2)
# because this runtime's Font property acts as a persistent "typing
# default" that otherwise leaks forward into text typed later.
$paraInfo = New-Object System.Collections.ArrayList

function Add-Paragraph([string]$text, [int]$boldChars) {
    $sel.TypeParagraph()
    $start = $d.Paragraphs.Last.Range.Start
    $sel.TypeText($text)
    if ($boldChars -gt 0) {
        [void]$paraInfo.Add(@{ Start = $start; BoldChars = $boldChars })
    }
}

Add-Paragraph "Lindsay Rizzardi Johns Hopkins Center for Epigenetics" 16
Add-Paragraph "Research Interest: Epigenetics and genome-wide data analysis" 0
Add-Paragraph "Goal: Learn how to manipulate, analyze, and interpret genome-wide data - - -" 0
Add-Paragraph "Girish Nadkarni; Institution: Mount Sinai Icahn School of Medicine" 15
Add-Paragraph "What I get out of this course: High-level idea of bioinformatics; data management for utilization in clinical and genomic research - - -" 0
Add-Paragraph "Payal Khincha" 13
Add-Paragraph "Pediatric Hematologist-Oncologist" 0
Add-Paragraph "Clinical Fellow, Clinical Genetics Branch, NCI" 0
Add-Paragraph "I am hoping to get a good understanding of basic bioinformatics so I can smartly handle exome sequencing data in my upcoming projects! - - -" 0

# Phase 2: now that no more typing will happen, apply bold formatting to
# the recorded leading-name ranges.
foreach ($info in $paraInfo) {
    $boldRange = $d.Range($info.Start, $info.Start + $info.BoldChars)
    $boldRange.Font.Bold = 1
}

Write-Host "Appended student participant paragraphs."
